# Fruta / hortaliza, semanal
# Insert two new weekly records for "Vega Monumental Concepción - Limón"
# right before the existing row 464, shifting the remaining rows down by
# two (old row 464 -> new row 466, ... old row 483 -> new row 485).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 464 (pushes everything below down by 2)
$ws.Rows.Item(464).Insert()
$ws.Rows.Item(464).Insert()

# Copy the static / repeated columns (and styles) from the row that is now
# just above (row 463) so the new rows keep the same Mercado / Producto
# metadata and number formats. Done one destination row at a time because
# pasting a single source row into a multi-row destination does not tile.
$ws.Range("A463:T463").Copy()
$ws.Range("A464:T464").PasteSpecial()
$ws.Range("A463:T463").Copy()
$ws.Range("A465:T465").PasteSpecial()

$newDate = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0

# New row 464: 1a amarillo
$ws.Cells.Item(464, 4).Value = $newDate
$ws.Cells.Item(464, 12).Value = "1a amarillo"
$ws.Cells.Item(464, 13).Value = 600
$ws.Cells.Item(464, 14).Value = 12000
$ws.Cells.Item(464, 15).Value = 13000
$ws.Cells.Item(464, 16).Value = 12500
$ws.Cells.Item(464, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(464, 19).Value = 781

# New row 465: 2a amarillo
$ws.Cells.Item(465, 4).Value = $newDate
$ws.Cells.Item(465, 12).Value = "2a amarillo"
$ws.Cells.Item(465, 13).Value = 300
$ws.Cells.Item(465, 14).Value = 10000
$ws.Cells.Item(465, 15).Value = 10000
$ws.Cells.Item(465, 16).Value = 10000
$ws.Cells.Item(465, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(465, 19).Value = 625
